$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 5
$ws.Range('B5').Value = '1: Title=Connect Plus; MaxDataLimit=50
2: Title=Connect Next First Responder; MaxDataLimit=50
3: Title=Connect Next 55; MaxDataLimit=50
4: Title=Connect Plus First Responder; MaxDataLimit=50
5: Title=Connect Plus Military; MaxDataLimit=50
6: Title=Connect Next; MaxDataLimit=50'
$ws.Range('C5').Value = 'plan_mobile_hotspot_data(Title, DataLimit) :-
    plan(Title, _, _, features(Features)),
    member(feature(mobile_hotspot(_, high_speed_data_limit_gb(DataLimit)), _), Features).
most_mobile_hotspot_plan(Title, MaxDataLimit) :-
    findall(D, plan_mobile_hotspot_data(_, D), DataList),
    max_list(DataList, MaxDataLimit),
    plan_mobile_hotspot_data(Title, MaxDataLimit).'
$ws.Range('D5').Value = 'most_mobile_hotspot_plan(Title, MaxDataLimit).'
$ws.Rows.Item(5).RowHeight = 112

# Row 6
$ws.Range('A6').Value = 'Are there any plans for first responders? If so what are their prices for 4 lines.'
$ws.Range('B6').Value = '1: Title=Connect Next First Responder; Price=200
2: Title=Connect Plus First Responder; Price=160
3: Title=Connect First Responder; Price=120'
$ws.Range('C6').Value = 'plan_for_first_responder(Title, Price) :-
    plan(Title, category(first_responder), lines(Lines), _),
    member(line(4, monthly_price(Price), _), Lines).'
$ws.Range('D6').Value = 'plan_for_first_responder(Title, Price).'
$ws.Rows.Item(6).RowHeight = 64

# Row 7
$ws.Range('A7').Value = 'Find all available plans that can support 3 lines. Give their names and cost for 3 lines.'
$ws.Range('B7').Value = '1: Title=Connect Plus; Price=150
2: Title=Connect Next First Responder; Price=165
3: Title=Core Saver; Price=100
4: Title=Connect Next 55; Price=195
5: Title=Core; Price=90
6: Title=Connect Next Military; Price=165
7: Title=Connect Plus 55; Price=165
8: Title=Connect Plus First Responder; Price=135
9: Title=Connect Plus Military; Price=135
10: Title=Connect 55; Price=135
11: Title=Connect First Responder; Price=105
12: Title=Connect Next; Price=180
13: Title=Connect Military; Price=105'
$ws.Range('C7').Value = 'plan_for_lines(Num, Title, Price) :-
    plan(Title, _, lines(Ls), _),
    member(line(Num, monthly_price(Price), _), Ls).'
$ws.Range('D7').Value = 'plan_for_lines(3, Title, Price).'
$ws.Rows.Item(7).RowHeight = 208

# Row 8
$ws.Range('A8').Value = 'Are there any special plans for veterans and if so  list them and also list if taxes and fees are included in the price.'
$ws.Range('B8').Value = '1: Title=Connect Next Military; TaxesIncluded=yes
2: Title=Connect Plus Military; TaxesIncluded=yes
3: Title=Connect Military; TaxesIncluded=yes'
$ws.Range('C8').Value = 'veteran_plan(Title, TaxesIncluded) :- 
    plan(Title, category(military_veteran), _Lines, features(Features)), 
    member(feature(taxes_and_fees(_, included_in_monthly_price(TaxesIncluded)), _), Features).'
$ws.Range('D8').Value = 'veteran_plan(Title, TaxesIncluded).'
$ws.Rows.Item(8).RowHeight = 96

# Row 9
$ws.Range('A9').Value = 'What categories of plans are available?'
$ws.Range('B9').Value = 'Categories=[55_plus,all,first_responder,military_veteran]"'
$ws.Range('D9').Value = 'setof(Category, Title^Lines^Features^(plan(Title, category(Category), Lines, Features)), Categories).'
$ws.Rows.Item(9).RowHeight = 32

# Row 10
$ws.Range('B10').Value = 'Included=no'
$ws.Range('C10').Value = '%% Helper predicate to check if a feature applies to a given line number
feature_applies_for_line(applies_to_lines(all), _Line).
feature_applies_for_line(applies_to_lines(lines(Low,High)), Line) :- Line >= Low, Line =< High.
%% Predicate to check if the netflix feature is included for a given plan title and number of lines
plan_netflix_included(PlanTitle, NumLines, Included) :-
    plan(PlanTitle, _, lines(Lines), features(Features)),
    %% Ensure the plan supports the requested number of lines
    member(line(NumLines, monthly_price(_), _), Lines),
    (   %% If there''s a netflix feature that applies to the given line
        member(feature(netflix(_, included(Inc)), Appl), Features),
        feature_applies_for_line(Appl, NumLines)
    ->  Included = Inc
    ;   %% Otherwise, netflix is not included
        Included = no
    ).'
$ws.Range('D10').Value = 'plan_netflix_included(''Core'', 1, Included).'
$ws.Range('A10').Value = 'I want to purchase 1 line of ''Core''. Is netflix included?'
$ws.Rows.Item(10).RowHeight = 409.5

# Row 11
$ws.Range('B11').Value = 'Title=Connect 55; Lines=1; Price=60'
$ws.Range('C11').Value = 'candidate_55_netflix_plan(Title, Lines, Price) :-
    plan(Title, category("55_plus"), lines(LinesList), features(Features)),
    member(line(Lines, monthly_price(Price), _), LinesList),
    member(feature(netflix(_, included(yes)), _), Features).
cheapest_55_netflix_plan(Title, Lines, Price) :-
    setof((P, N, T), candidate_55_netflix_plan(T, N, P), Sorted),
    Sorted = [(Price, Lines, Title)|_].'
$ws.Range('D11').Value = 'cheapest_55_netflix_plan(Title, Lines, Price).'
$ws.Range('A11').Value = 'I am over 55 and am looking for the cheapest plan where netflix is included. Give the price and the number of lines'
$ws.Rows.Item(11).RowHeight = 192

[void]$ws.Range('A11').Select()
